$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert a new blank row at the very top. Everything that used to
#    be on rows 1-6 shifts down to rows 2-7 (and the old A1:D1 merge
#    moves with it, becoming A2:D2).
# ------------------------------------------------------------------
$ws.Rows.Item(1).Insert()

# ------------------------------------------------------------------
# 2) Fill in the three new requirement rows (rows 5-7, which used to
#    be the trailing blank placeholder rows of the original sheet).
# ------------------------------------------------------------------
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "Số lượng user"
$ws.Range("C5").Value = "30 users truy cập cùng lúc."
$ws.Range("D5").Value = ""

$ws.Range("A6").Value = 3
$ws.Range("B6").Value = "Tránh mất dữ liệu"
$ws.Range("C6").Value = "Tự động backup dữ liệu mỗi 12 tiếng."
$ws.Range("D6").Value = ""

$ws.Range("A7").Value = 4
$ws.Range("B7").Value = "Thời gian cài đặt"
$ws.Range("C7").Value = "Thời gian cài đặt không quá 2 tiếng."
$ws.Range("D7").Value = ""
$ws.Range("D4").Value = ""

# ------------------------------------------------------------------
# 3) Merge cells first (must happen before border formatting so that
#    every member cell of the merged area keeps its own complete
#    border instead of Excel normalizing it down to an outer box).
# ------------------------------------------------------------------
$ws.Range("A2:D2").Merge()
$ws.Range("A1:D1").Merge()

# ------------------------------------------------------------------
# 4) Row heights.
# ------------------------------------------------------------------
$ws.Range("A1:D7").RowHeight = 20.1

# ------------------------------------------------------------------
# 5) Column widths.
# ------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 71.6

# ------------------------------------------------------------------
# 6) Formatting for the whole table body (title + header + data,
#    rows 2-7): Helvetica 13pt, left/top aligned, full thin box
#    border on every individual cell.
# ------------------------------------------------------------------
$body = $ws.Range("A2:D7")
$body.Font.Name = "Helvetica"
$body.Font.Size = 13
$body.VerticalAlignment = -4160    # xlTop
$body.HorizontalAlignment = -4131  # xlLeft
$body.Borders.LineStyle = 1
$body.Borders.Weight = 2
$body.Borders.Item(11).LineStyle = 1   # xlInsideVertical
$body.Borders.Item(11).Weight = 2
$body.Borders.Item(12).LineStyle = 1   # xlInsideHorizontal
$body.Borders.Item(12).Weight = 2

# Bold the title cell (row 2, A2 only - the rest of the merged row
# stays regular weight) and the whole header row (row 3).
$ws.Range("A2").Font.Bold = $true
$ws.Range("A3:D3").Font.Bold = $true

# ------------------------------------------------------------------
# 7) Row 1 formatting: Helvetica 13pt, centered, top aligned, bottom
#    border only (no box).
# ------------------------------------------------------------------
$row1 = $ws.Range("A1:D1")
$row1.Font.Name = "Helvetica"
$row1.Font.Size = 13
$row1.VerticalAlignment = -4160    # xlTop
$row1.HorizontalAlignment = -4108  # xlCenter
$row1.Borders.Item(9).LineStyle = 1    # xlEdgeBottom
$row1.Borders.Item(9).Weight = 2

# ------------------------------------------------------------------
# 8) Selection / active cell, matching the saved view state.
# ------------------------------------------------------------------
$ws.Range("I13").Select()
